$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the title
#    (Heading1) paragraph at the top of the document.
# ---------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs.Item(2)
$metaPara.Style = "Normal"

$metaStart = $metaPara.Range.Start
$metaFull = "Meta description: Read our review of the Book of Dead slot game. Play for free with a wide range of betting options and unique features like wild and scatter symbols represented by the same symbol."

# Exclude the trailing paragraph mark from the text range we set so the
# paragraph mark itself keeps its own (default) formatting.
$metaTextRange = $d.Range($metaStart, $metaPara.Range.End - 1)
$metaTextRange.Text = $metaFull

# Bold just the "Meta description" label, leaving the rest of the
# sentence (starting at the colon) in regular formatting.
$labelLen = ("Meta description").Length
$labelRange = $d.Range($metaStart, $metaStart + $labelLen)
$labelRange.Bold = 1

# ---------------------------------------------------------------------
# 2) Locate and remove the old duplicate bold "Play Book of Dead for
#    Free - Slot Game Review" paragraph further down the document
#    (the Heading1 at the very top must stay untouched).
# ---------------------------------------------------------------------
$titleText = "Play Book of Dead for Free - Slot Game Review"

$scan = $d.Content
$scan.Start = $titlePara.Range.End
$scan.End = $d.Content.End
$found = $scan.Find.Execute($titleText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # Expand the found range to the whole paragraph (including its
    # paragraph mark) so the entire paragraph is removed cleanly.
    $dupRange = $d.Range($scan.Start, $scan.End)
    $dupRange.Expand(4) | Out-Null
    $dupRange.Delete()
}

# ---------------------------------------------------------------------
# 3) Replace the text of the final (italic) paragraph with the new
#    image-generation prompt, keeping its italic formatting intact.
# ---------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastRange = $d.Range($lastPara.Range.Start, $lastPara.Range.End - 1)
$newPrompt = 'Create a cartoon-style feature image for the game "Book of Dead" featuring a happy Maya warrior with glasses. The image should have vibrant colors and showcase the ancient Egyptian theme of the game. The Maya warrior should be holding the Book of Dead and standing in front of the pyramids. The background should have a sunset hue with Egyptian hieroglyphics in the sky. Make sure to highlight the warrior''s glasses which should be oversized and reflective. The overall style of the image should be fun and eye-catching, inviting players to dive into the adventure-packed world of "Book of Dead."'
$lastRange.Text = $newPrompt
